$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a daily price log for "Femacal de La Calera - Piña".
# Two new observation rows (week of 2021-11-09, serial 44509) are inserted
# right before the current row 369, pushing the existing rows 369-382 down
# to 371-384.

$ws.Rows.Item(369).Insert()
$ws.Rows.Item(369).Insert()

# Match the date-column number format used by the surrounding rows (now at 371/372).
$dateFmt = $ws.Cells.Item(371, 4).NumberFormat
$ws.Cells.Item(369, 4).NumberFormat = $dateFmt
$ws.Cells.Item(370, 4).NumberFormat = $dateFmt

# --- Row 369: Primera, $/caja 12 unidades ---
$ws.Cells.Item(369, 1).Value  = 3
$ws.Cells.Item(369, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(369, 3).Value  = "Coquimbo"
$ws.Cells.Item(369, 4).Value  = 44509
$ws.Cells.Item(369, 5).Value  = 5
$ws.Cells.Item(369, 6).Value  = "Fruta"
$ws.Cells.Item(369, 7).Value  = 100108
$ws.Cells.Item(369, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(369, 9).Value  = 100108005
$ws.Cells.Item(369, 10).Value = "Piña"
$ws.Cells.Item(369, 11).Value = "Caramelo"
$ws.Cells.Item(369, 12).Value = "Primera"
$ws.Cells.Item(369, 13).Value = 216
$ws.Cells.Item(369, 14).Value = 20000
$ws.Cells.Item(369, 15).Value = 20000
$ws.Cells.Item(369, 16).Value = 20000
$ws.Cells.Item(369, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(369, 18).Value = "Ecuador"
$ws.Cells.Item(369, 19).Value = 1667
$ws.Cells.Item(369, 20).Value = 12

# --- Row 370: Segunda, $/caja 14 unidades ---
$ws.Cells.Item(370, 1).Value  = 3
$ws.Cells.Item(370, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(370, 3).Value  = "Coquimbo"
$ws.Cells.Item(370, 4).Value  = 44509
$ws.Cells.Item(370, 5).Value  = 5
$ws.Cells.Item(370, 6).Value  = "Fruta"
$ws.Cells.Item(370, 7).Value  = 100108
$ws.Cells.Item(370, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(370, 9).Value  = 100108005
$ws.Cells.Item(370, 10).Value = "Piña"
$ws.Cells.Item(370, 11).Value = "Caramelo"
$ws.Cells.Item(370, 12).Value = "Segunda"
$ws.Cells.Item(370, 13).Value = 108
$ws.Cells.Item(370, 14).Value = 20000
$ws.Cells.Item(370, 15).Value = 20000
$ws.Cells.Item(370, 16).Value = 20000
$ws.Cells.Item(370, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(370, 18).Value = "Ecuador"
$ws.Cells.Item(370, 19).Value = 1429
$ws.Cells.Item(370, 20).Value = 14
